$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19 ---------------------------------------------------------------
$ws.Range("A19").Value = "opals"

$modelCode = "model = keras.Sequential()`ninput_shape = (1, 220500)`nsr = 22050  # Sampling rate of the audio`nmodel.add(Melspectrogram(n_dft=512, n_hop=256, input_shape=input_shape,`n                         padding='same', sr=sr, n_mels=128,`n                         fmin=0.0, fmax=sr/2, power_melgram=1.0,`n                         return_decibel_melgram=False, trainable_fb=False,`n                         trainable_kernel=False,`n                         name='trainable_stft'))`nmodel.add(keras.layers.Conv2D(128, [7,11], strides=[2,2], padding = 'SAME'))`nmodel.add(keras.layers.LeakyReLU(alpha=0.1))`nmodel.add(keras.layers.MaxPool2D(pool_size=(2,2), padding='SAME'))`nmodel.add(keras.layers.Conv2D(128, [7,11], strides=[2,2], padding = 'SAME'))`nmodel.add(keras.layers.LeakyReLU(alpha=0.1))`n# model.add(keras.layers.Dropout(0.5))`nmodel.add(keras.layers.Flatten())`nmodel.add(keras.layers.Dense(1, activation='sigmoid'))`n"
$ws.Range("C19").Value = $modelCode
# the long pasted-in code block carries its own (no-op) alignment style,
# distinct from the plain default style used elsewhere in the sheet
$ws.Range("C19").VerticalAlignment = -4107
# undo the implicit row-height auto-expand caused by the multi-line paste
$ws.Rows.Item(19).AutoFit()

$ws.Range("E19").Value = 1.33
$ws.Range("F19").Value = 0.81
$ws.Range("G19").Value = 46
$ws.Range("H19").Value = 1
$ws.Range("I19").Value = 16
$ws.Range("J19").Value = 26
$ws.Range("K19").NumberFormat = $ws.Range("K2").NumberFormat
$ws.Range("K19").Formula = "=(G19+J19)/SUM(G19:J19)"

$ws.Range("L19").Value = "CORRECTED CONFUSION MATRIX AT THIS POINT"

# --- Row 20 -----------------------------------------------------------
$ws.Range("A20").Value = "pringle"
$ws.Range("C20").Value = "same as above but with early stopping"
$ws.Range("E20").Value = 0.53
$ws.Range("F20").Value = 0.78
$ws.Range("G20").Value = 45
$ws.Range("H20").Value = 2
$ws.Range("I20").Value = 18
$ws.Range("J20").Value = 24
$ws.Range("K20").NumberFormat = $ws.Range("K2").NumberFormat
$ws.Range("K20").Formula = "=(G20+J20)/SUM(G20:J20)"

# --- Row 21 -----------------------------------------------------------
$ws.Range("C21").Value = "same as above with a 0.5 dropout at the end"
$ws.Range("E21").Value = 4.45
$ws.Range("F21").Value = 0.68
$ws.Range("G21").Value = 33
$ws.Range("H21").Value = 14
$ws.Range("I21").Value = 14
$ws.Range("J21").Value = 28
$ws.Range("K21").NumberFormat = $ws.Range("K2").NumberFormat
$ws.Range("K21").Formula = "=(G21+J21)/SUM(G21:J21)"

# --- Row 22 -----------------------------------------------------------
$ws.Range("C22").Value = "same as above bu 0.2 dropout at end"
$ws.Range("A22").Value = "Query"
$ws.Range("E22").Value = 0.58
$ws.Range("F22").Value = 0.84
$ws.Range("G22").Value = 44
$ws.Range("H22").Value = 3
$ws.Range("I22").Value = 11
$ws.Range("J22").Value = 31
$ws.Range("K22").NumberFormat = $ws.Range("K2").NumberFormat
$ws.Range("K22").Formula = "=(G22+J22)/SUM(G22:J22)"

# --- Row 23 -----------------------------------------------------------
$ws.Range("C23").Value = "same as before but 3nd conv kernel size 64, no dropiut"

# --- Extend the conditional formatting ranges down through row 22 ---------
$ws.Range("K2:K22").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("K2:K22"))
$ws.Range("F2:F22").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("F2:F22"))
$ws.Range("E2:E22").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("E2:E22"))

# --- Match the saved selection position ------------------------------------
$null = $ws.Range("C24").Select()
